# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period table (rows 16-23, columns C/D/E) is re-sequenced:
# instead of alternating period 1801/1712 per worker, the rows are now
# grouped by period (all workers for period 1712 first, then the same
# workers for period 1801).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Doc number (C) / Name (D) for each of the 4 workers, in the order they
# must now appear for each period block.
$docs = @("23071486", "32941250", "1049533950", "1049532082")
$names = @(
    "DANILZA MARIA MERIÑO DIAZ",
    "LUCILA HERRERA GUZMAN",
    "ROSANA MARIA PADILLA JULIO",
    "LORENA MARIA BOLAÑOS UTRIA"
)

# First block: rows 16-19 -> periodo 1712 for each of the 4 workers
$row = 16
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($row, 3).Value = $docs[$i]
    $ws.Cells.Item($row, 4).Value = $names[$i]
    $ws.Cells.Item($row, 5).Value = "1712"
    $row++
}

# Second block: rows 20-23 -> periodo 1801 for each of the 4 workers
for ($i = 0; $i -lt 4; $i++) {
    $ws.Cells.Item($row, 3).Value = $docs[$i]
    $ws.Cells.Item($row, 4).Value = $names[$i]
    $ws.Cells.Item($row, 5).Value = "1801"
    $row++
}
